# Append two new paragraphs ("Chapter 12, Question 2" heading + its answer
# write-up) to the end of the document body, matching the formatting of the
# existing "Chapter N, Question M" sections already in the document.
$d = $word.ActiveDocument

# Position right at the end of the body content (after the last paragraph's
# text, before the section break) so InsertXML appends rather than replaces.
$endPos = $d.Content.End
$target = $d.Range($endPos, $endPos)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Chapter 12, Question 2</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/><w:t xml:space="preserve">This question asked me to add a log statement at the beginning of the runner’s </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>stopFalling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> method that prints the vertical velocity, in pixels per second. I added this log statement to the beginning of the method. It is triggered as soon as the runner stops falling and lands on something. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>The velocity greatly depends on where the runner is jumping from and where she’s falling to. Jumping from the top track to the bottom track, for example, yields a velocity of ~460 pixels per second. Falling from the top track to the bottom track yields a velocity of ~348 pixels per second. Falling from the middle track to the bottom track yields a velocity of ~244 pixels per second. It varies depending on the time spent falling.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)
